$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

# The old D column held a buggy INDEX/MATCH lookup formula trying to resolve
# the column name; the correct calculation already existed one column over
# in column E. Move that correct formula into D, then drop the now-redundant
# column E (which also removes the stray "col_5" header in E1).
$ws3.Range("D2").Formula = $ws3.Range("E2").Formula
$ws3.Range("D3").Formula = $ws3.Range("E3").Formula
$ws3.Columns.Item(5).Delete()

$ws3.Activate()
$ws3.Range("D2").Select()
